# Smartwatch EXTRA-BOTTOM BOM update:
#  - Swap the right-angle rotary encoder (PEC09-2220F-S0012 / Bourns) for a
#    vertical rotary encoder (EVQ-VUA00112B / Panasonic) on row 25 (ENC1).
#  - Add a new tactile button (B1, C&K KT11P2SA1M35LFG) as a new BOM line
#    right after the encoder.
#  - Keep the trailing "Mouser" divider + AD5171 rows intact (they just
#    shift down/up by the net row-count change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the currency number format on the Price-per-unit column for
#    the rows we're touching so K24:K26 match the rest of the table
#    (Accounting/Currency format instead of plain General).
# ---------------------------------------------------------------------
$currencyFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

# ---------------------------------------------------------------------
# 2. Row 25 (ENC1): replace the right-angle encoder with the vertical one.
# ---------------------------------------------------------------------
$ws.Range("B25").Value = "Encoder"
$ws.Range("C25").Value = "ENC1"
$ws.Range("D25").Value = "Rotary Encoder Mechanical 12 Quadrature (Incremental) Vertical"
$ws.Range("E25").Value = "THT"
$ws.Range("F25").Value = ""
$ws.Range("G25").Value = "Panasonic Electronic Components"
$ws.Range("H25").Value = "EVQ-VUA00112B"
$ws.Range("I25").Value = "Digi-Key"
$ws.Range("J25").Value = "P123413-ND"
$ws.Range("K25").Value = 1.18
$ws.Range("L25").Value = 1
$ws.Range("M25").Formula = "=L25*K25"
$ws.Range("N25").Value = 0.80625
$ws.Range("O25").Formula = "=N25*L25"

$ws.Range("K25").NumberFormat = $currencyFormat
$ws.Range("L25").NumberFormat = "General"

# ---------------------------------------------------------------------
# 3. Insert a brand-new row 26 for the tactile button (B1) right after
#    the encoder row; Excel re-targets every downstream formula/shared
#    range automatically.
# ---------------------------------------------------------------------
$ws.Rows(26).Insert()

$ws.Range("A26").Formula = "=A25+1"
$ws.Range("B26").Value = "Tactile Button"
$ws.Range("C26").Value = "B1"
$ws.Range("D26").Value = "SWITCH TACTILE SPST-NO 1VA 32V"
$ws.Range("E26").Value = "SMD"
$ws.Range("G26").Value = "C&K"
$ws.Range("H26").Value = "KT11P2SA1M35LFG"
$ws.Range("I26").Value = "Digi-Key"
$ws.Range("J26").Value = "CKN1843CT-ND"
$ws.Range("K26").Value = 4.47
$ws.Range("L26").Value = 1
$ws.Range("M26").Formula = "=L26*K26"
$ws.Range("N26").Value = 3.3904
$ws.Range("O26").Formula = "=N26*L26"

$ws.Range("A26").Style = $ws.Range("A25").Style
$ws.Range("B26:E26").Style = $ws.Range("B25").Style
$ws.Range("G26:H26").Style = $ws.Range("G25").Style
$ws.Range("I26").Style = $ws.Range("I25").Style
$ws.Range("J26").Style = $ws.Range("J25").Style
$ws.Range("M26").Style = $ws.Range("M25").Style
$ws.Range("N26").Style = $ws.Range("N25").Style
$ws.Range("O26").Style = $ws.Range("O25").Style

$ws.Range("K26").NumberFormat = $currencyFormat
$ws.Range("L26").NumberFormat = "General"

# ---------------------------------------------------------------------
# 4. The two now-blank spacer rows that used to sit between the encoder
#    and the "Mouser" divider (old rows 26 & 27) are no longer needed --
#    delete them so the "Mouser"/AD5171 block shifts back up by one,
#    giving a net -1 row overall (29 -> 28).
# ---------------------------------------------------------------------
$ws.Rows("27:28").Delete()

# ---------------------------------------------------------------------
# 5. Match the cursor position recorded in the saved file.
# ---------------------------------------------------------------------
$ws.Range("C14").Select()
